# DAX and EP Global Objects
#
# The "RVL" sheet's Global Objects map had four now-unused Range
# parameters (fromRow/fromCol/toRow/toCol) removed, and the remaining
# "Functions" actions for launching/selecting the DAX client were
# renamed to a dedicated "DAX" object (DaxLaunch -> Launch,
# DaxChangeCompany -> ChangeCompany, DaxOpenModule -> OpenModule,
# DaxNavigate -> Navigate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# Remove the fromRow/fromCol/toRow/toCol Param rows of the "Map Range"
# block; this shifts every following row up by four.
$ws.Rows("9:12").Delete()

# Rename the "Functions" object to "DAX" and drop the "Dax" prefix from
# the action names for the four rows that used to launch/select the
# Dax client (now rows 11-14 after the delete above).
$ws.Range("C11").Value = "DAX"
$ws.Range("D11").Value = "Launch"

$ws.Range("C12").Value = "DAX"
$ws.Range("D12").Value = "ChangeCompany"

$ws.Range("C13").Value = "DAX"
$ws.Range("D13").Value = "OpenModule"

$ws.Range("C14").Value = "DAX"
$ws.Range("D14").Value = "Navigate"
